# BENTO_GONCALVES.xlsx - automatic update
#
# 1. "Paineis DARQ" -> "PAINEIS DARQ"
# 2. "Recolhimento x Eliminacao" -> "RECOLHIMENTO X ELIMINAÇÃO"
# 3. Remove the "Desarquivamentos Pendentes" sheet entirely

$wb = $excel.ActiveWorkbook

# 1 & 2: rename the two sheets (case / accent normalisation)
$wb.Worksheets.Item("Paineis DARQ").Name = "PAINEIS DARQ"
$wb.Worksheets.Item("Recolhimento x Eliminacao").Name = "RECOLHIMENTO X ELIMINAÇÃO"

# 3: delete the "Desarquivamentos Pendentes" worksheet (suppress the
# "this will permanently delete" confirmation prompt around the call)
$excel.DisplayAlerts = $false
$wb.Worksheets.Item("Desarquivamentos Pendentes").Delete()
$excel.DisplayAlerts = $true
